# Auto-generated COM-interop script: add 2022-Q1 sheet + update 总计 summary sheet
$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" worksheet, positioned after "2021-Q4" ---
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# Copy the header row (with its style) from the "2021-Q4" sheet as a formatting template
$afterSheet.Range("A1:H1").Copy($q1.Range("A1:H1"))
$afterSheet.Range("A2").Copy($q1.Range("A2:A17"))
$afterSheet.Range("H2").Copy($q1.Range("H2:H17"))

# Header text
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Force columns B:G (fund code/name/size/position/ratio/value) to be stored as literal text,
# matching the source data (these look numeric but are text in the workbook).
$q1.Range("B2:G17").NumberFormat = "@"

# Data rows 2..17
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "320003"
$q1.Range("C2").Value = "诺安先锋混合"
$q1.Range("D2").Value = "45.79"
$q1.Range("E2").Value = "69.96"
$q1.Range("F2").Value = "6.82"
$q1.Range("G2").Value = "3.1229"
$q1.Range("H2").Value = 1
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "000362"
$q1.Range("C3").Value = "国泰聚信价值优势灵活配置混合A"
$q1.Range("D3").Value = "56.15"
$q1.Range("E3").Value = "89.00"
$q1.Range("F3").Value = "4.69"
$q1.Range("G3").Value = "2.6334"
$q1.Range("H3").Value = 6
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "160106"
$q1.Range("C4").Value = "南方高增长混合(LOF)"
$q1.Range("D4").Value = "20.27"
$q1.Range("E4").Value = "87.01"
$q1.Range("F4").Value = "4.88"
$q1.Range("G4").Value = "0.9892"
$q1.Range("H4").Value = 6
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "000363"
$q1.Range("C5").Value = "国泰聚信价值优势灵活配置混合C"
$q1.Range("D5").Value = "17.09"
$q1.Range("E5").Value = "89.00"
$q1.Range("F5").Value = "4.69"
$q1.Range("G5").Value = "0.8015"
$q1.Range("H5").Value = 6
$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "020010"
$q1.Range("C6").Value = "国泰金牛创新混合"
$q1.Range("D6").Value = "16.99"
$q1.Range("E6").Value = "84.21"
$q1.Range("F6").Value = "4.07"
$q1.Range("G6").Value = "0.6915"
$q1.Range("H6").Value = 7
$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "012173"
$q1.Range("C7").Value = "国泰兴泽优选一年持有期混合A"
$q1.Range("D7").Value = "16.89"
$q1.Range("E7").Value = "89.83"
$q1.Range("F7").Value = "3.25"
$q1.Range("G7").Value = "0.5489"
$q1.Range("H7").Value = 10
$q1.Range("A8").Value = 6
$q1.Range("B8").Value = "001743"
$q1.Range("C8").Value = "诺安优选回报灵活配置混合"
$q1.Range("D8").Value = "6.13"
$q1.Range("E8").Value = "71.32"
$q1.Range("F8").Value = "7.87"
$q1.Range("G8").Value = "0.4824"
$q1.Range("H8").Value = 2
$q1.Range("A9").Value = 7
$q1.Range("B9").Value = "160105"
$q1.Range("C9").Value = "南方积极配置混合(LOF)"
$q1.Range("D9").Value = "7.84"
$q1.Range("E9").Value = "88.09"
$q1.Range("F9").Value = "5.34"
$q1.Range("G9").Value = "0.4187"
$q1.Range("H9").Value = 5
$q1.Range("A10").Value = 8
$q1.Range("B10").Value = "005244"
$q1.Range("C10").Value = "国泰聚优价值灵活配置混合A"
$q1.Range("D10").Value = "7.72"
$q1.Range("E10").Value = "83.97"
$q1.Range("F10").Value = "3.32"
$q1.Range("G10").Value = "0.2563"
$q1.Range("H10").Value = 5
$q1.Range("A11").Value = 9
$q1.Range("B11").Value = "012174"
$q1.Range("C11").Value = "国泰兴泽优选一年持有期混合C"
$q1.Range("D11").Value = "7.14"
$q1.Range("E11").Value = "89.83"
$q1.Range("F11").Value = "3.25"
$q1.Range("G11").Value = "0.2320"
$q1.Range("H11").Value = 10
$q1.Range("A12").Value = 10
$q1.Range("B12").Value = "008185"
$q1.Range("C12").Value = "诺安研究优选混合"
$q1.Range("D12").Value = "2.59"
$q1.Range("E12").Value = "94.06"
$q1.Range("F12").Value = "6.24"
$q1.Range("G12").Value = "0.1616"
$q1.Range("H12").Value = 5
$q1.Range("A13").Value = 11
$q1.Range("B13").Value = "005245"
$q1.Range("C13").Value = "国泰聚优价值灵活配置混合C"
$q1.Range("D13").Value = "4.52"
$q1.Range("E13").Value = "83.97"
$q1.Range("F13").Value = "3.32"
$q1.Range("G13").Value = "0.1501"
$q1.Range("H13").Value = 5
$q1.Range("A14").Value = 12
$q1.Range("B14").Value = "003131"
$q1.Range("C14").Value = "国寿安保强国智造灵活配置混合"
$q1.Range("D14").Value = "5.76"
$q1.Range("E14").Value = "86.19"
$q1.Range("F14").Value = "2.20"
$q1.Range("G14").Value = "0.1267"
$q1.Range("H14").Value = 5
$q1.Range("A15").Value = 13
$q1.Range("B15").Value = "005683"
$q1.Range("C15").Value = "国寿安保华兴灵活配置混合"
$q1.Range("D15").Value = "3.43"
$q1.Range("E15").Value = "89.19"
$q1.Range("F15").Value = "2.03"
$q1.Range("G15").Value = "0.0696"
$q1.Range("H15").Value = 8
$q1.Range("A16").Value = 14
$q1.Range("B16").Value = "000554"
$q1.Range("C16").Value = "南方中国梦灵活配置混合"
$q1.Range("D16").Value = "1.39"
$q1.Range("E16").Value = "88.18"
$q1.Range("F16").Value = "4.97"
$q1.Range("G16").Value = "0.0691"
$q1.Range("H16").Value = 8
$q1.Range("A17").Value = 15
$q1.Range("B17").Value = "001744"
$q1.Range("C17").Value = "诺安进取回报灵活配置混合"
$q1.Range("D17").Value = "0.04"
$q1.Range("E17").Value = "62.10"
$q1.Range("F17").Value = "5.80"
$q1.Range("G17").Value = "0.0023"
$q1.Range("H17").Value = 2

# --- Step 2: update the "总计" summary sheet: insert a new top data row for 2022-Q1 ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

# Rewrite every data row (index + date + count + market value) top to bottom
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 16
$total.Range("D2").Value = 10.76
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 23
$total.Range("D3").Value = 12.69
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 3.16
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 2.01
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 3
$total.Range("D6").Value = 2.03
$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 6
$total.Range("D7").Value = 2.3

Write-Host "edit complete"
